$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New raw data values in column H
$ws.Range("H4").Value2 = 0
$ws.Range("H5").Value2 = 15
$ws.Range("H6").Value2 = 0
$ws.Range("H7").Value2 = 15
$ws.Range("H8").Value2 = 0
$ws.Range("H9").Value2 = 15
$ws.Range("H10").Value2 = 0

# New running-total formulas in column I (bold, like columns C and F)
$ws.Range("I6").Formula = "=H4+H5+H6"
$ws.Range("I8").Formula = "=I6+H7+H8"
$ws.Range("I10").Formula = "=I8+H9+H10"

# Apply bold style to the running-total cells, matching style used by C6/F6/C8/F8/C10/F10
$ws.Range("I6").Font.Bold = $true
$ws.Range("I8").Font.Bold = $true
$ws.Range("I10").Font.Bold = $true

# Sum formula in H11 (bold, like B13/E13 totals but placed at row 11 per the diff)
$ws.Range("H11").Formula = "=SUM(H4:H10)"
$ws.Range("H11").Font.Bold = $true

# Update selection to match the recorded end-state of the edit
$ws.Range("I6").Select()
